$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (2-10) for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$data = @{
    2  = @{ E=3; G=76.08252666666667;  H=228.24758;         I=0.95878149807566;  J=0.95878149807566;  K=3; M=8.806900666666666;  N=26.420702;         O=0.1733678197953833; P=0.1733678197953834; Q=670.051254822351;  R=6030.461293401159;  S=0.1662218579815287;   T=0.1662218579815287 }
    3  = @{ E=3; G=76.08252666666667;  H=228.24758;         I=0.95878149807566;  J=0.95878149807566;  K=3; M=18.76689066666667;  N=56.30067200000001; O=0.3694347242421866; P=0.3694347242421866; Q=1427.832459597085; R=12850.49213637376;  S=0.354207178350092;    T=0.354207178350092 }
    4  = @{ E=3; G=76.08252666666667;  H=228.24758;         I=0.95878149807566;  J=0.95878149807566;  K=3; M=23.225144;           N=69.675432;         O=0.4571974559624301; P=0.4571974559624301; Q=1767.027637717173; R=15903.24873945456;  S=0.4383524617440394;   T=0.4383524617440394 }
    5  = @{ E=3; G=0.9347409999999999; H=2.804223;          I=0.01177947704364805; J=0.01177947704364805; K=3; M=8.806900666666666; N=26.420702;       O=0.1733678197953833; P=0.1733678197953834; Q=8.232171136060664; R=74.089540224546;    S=0.00204218225338703;  T=0.00204218225338703 }
    6  = @{ E=3; G=0.9347409999999999; H=2.804223;          I=0.01177947704364805; J=0.01177947704364805; K=3; M=18.76689066666667; N=56.30067200000001; O=0.3694347242421866; P=0.3694347242421866; Q=17.54218214865067; R=157.879639337856;   S=0.004351747853337284; T=0.004351747853337284 }
    7  = @{ E=3; G=0.9347409999999999; H=2.804223;          I=0.01177947704364805; J=0.01177947704364805; K=3; M=23.225144;        N=69.675432;        O=0.4571974559624301; P=0.4571974559624301; Q=21.709494327704;   R=195.385448949336;   S=0.005385546936923735; T=0.005385546936923735 }
    8  = @{ E=3; G=2.336085333333334;  H=7.008256;          I=0.02943902488069198; J=0.02943902488069198; K=3; M=8.806900666666666; N=26.420702;       O=0.1733678197953833; P=0.1733678197953834; Q=20.57367147952355; R=185.163043315712;   S=0.005103779560467613; T=0.005103779560467614 }
    9  = @{ E=3; G=2.336085333333334;  H=7.008256;          I=0.02943902488069198; J=0.02943902488069198; K=3; M=18.76689066666667; N=56.30067200000001; O=0.3694347242421866; P=0.3694347242421866; Q=43.84105803867023; R=394.5695223480321;  S=0.01087579803875731;  T=0.01087579803875731 }
    10 = @{ E=3; G=2.336085333333334;  H=7.008256;          I=0.02943902488069198; J=0.02943902488069198; K=3; M=23.225144;        N=69.675432;        O=0.4571974559624301; P=0.4571974559624301; Q=54.25591826295467; R=488.303264366592;   S=0.01345944728146706;  T=0.01345944728146706 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
